$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Options Sheet")
$ws = $wb.Worksheets.Item("My Template")

# PasteSpecial constants used below:
#   -4122 = xlPasteFormats (formats/styles only)
#   -4163 = xlPasteValues  (values only, no re-interpretation/auto date-parsing)

# --- Row 2: shift A:D -> B:E, then write new header key "processing" into A2 ---
$ws.Range("D2").Copy()
$ws.Range("E2").PasteSpecial(-4163)
$ws.Range("C2").Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("B2").Copy()
$ws.Range("C2").PasteSpecial(-4163)
$ws.Range("A2").Copy()
$ws.Range("B2").PasteSpecial(-4163)
$ws.Range("A2").Value = "processing"

# --- Row 3: shift A:D -> B:E (value + the blue header style), new A3 = "Processing" ---
$ws.Range("D3").Copy()
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("D3").Copy()
$ws.Range("E3").PasteSpecial(-4163)

$ws.Range("C3").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("D3").PasteSpecial(-4163)

$ws.Range("B3").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("C3").PasteSpecial(-4163)

$ws.Range("A3").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("A3").Copy()
$ws.Range("B3").PasteSpecial(-4163)

$ws.Range("B3").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Value = "Processing"

# --- Row 4: shift A:D -> B:E, new A4 = "Process1" ---
$ws.Range("D4").Copy()
$ws.Range("E4").PasteSpecial(-4163)
$ws.Range("C4").Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("B4").Copy()
$ws.Range("C4").PasteSpecial(-4163)
$ws.Range("A4").Copy()
$ws.Range("B4").PasteSpecial(-4163)
$ws.Range("A4").Value = "Process1"

# --- Row 5: shift A:D -> B:E, new A5 = "Process2" ---
$ws.Range("D5").Copy()
$ws.Range("E5").PasteSpecial(-4163)
$ws.Range("C5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("B5").Copy()
$ws.Range("C5").PasteSpecial(-4163)
$ws.Range("A5").Copy()
$ws.Range("B5").PasteSpecial(-4163)
$ws.Range("A5").Value = "Process2"

# --- Update the selection shown on "My Template" without leaving it as the active tab ---
$ws.Range("H1").Select()
$ws1.Activate()
